$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.983.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3794"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08315"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.106"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.37"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.195"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.42"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.169"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.007"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.024.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.263"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.567"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.081.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.38"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1054"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.035"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.581"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.595"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.628"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02414"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06508"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2150"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.205"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6367"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.230"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.852"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6065"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.285"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.659"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.986"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.64"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.09%  "
